$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell style (bold, bordered, centered, General format) onto A2:A39
# so the date column switches away from the custom date-time number format.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A39").PasteSpecial(-4122)

# Replace each date serial with its "<year>Q4" text label.
$ws.Range("A2").Value = "1987Q4"
$ws.Range("A3").Value = "1988Q4"
$ws.Range("A4").Value = "1989Q4"
$ws.Range("A5").Value = "1990Q4"
$ws.Range("A6").Value = "1991Q4"
$ws.Range("A7").Value = "1992Q4"
$ws.Range("A8").Value = "1993Q4"
$ws.Range("A9").Value = "1994Q4"
$ws.Range("A10").Value = "1995Q4"
$ws.Range("A11").Value = "1996Q4"
$ws.Range("A12").Value = "1997Q4"
$ws.Range("A13").Value = "1998Q4"
$ws.Range("A14").Value = "1999Q4"
$ws.Range("A15").Value = "2000Q4"
$ws.Range("A16").Value = "2001Q4"
$ws.Range("A17").Value = "2002Q4"
$ws.Range("A18").Value = "2003Q4"
$ws.Range("A19").Value = "2004Q4"
$ws.Range("A20").Value = "2005Q4"
$ws.Range("A21").Value = "2006Q4"
$ws.Range("A22").Value = "2007Q4"
$ws.Range("A23").Value = "2008Q4"
$ws.Range("A24").Value = "2009Q4"
$ws.Range("A25").Value = "2010Q4"
$ws.Range("A26").Value = "2011Q4"
$ws.Range("A27").Value = "2012Q4"
$ws.Range("A28").Value = "2013Q4"
$ws.Range("A29").Value = "2014Q4"
$ws.Range("A30").Value = "2015Q4"
$ws.Range("A31").Value = "2016Q4"
$ws.Range("A32").Value = "2017Q4"
$ws.Range("A33").Value = "2018Q4"
$ws.Range("A34").Value = "2019Q4"
$ws.Range("A35").Value = "2020Q4"
$ws.Range("A36").Value = "2021Q4"
$ws.Range("A37").Value = "2022Q4"
$ws.Range("A38").Value = "2023Q4"
$ws.Range("A39").Value = "2024Q4"

$excel.CutCopyMode = 0
